$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly update rotates the data rows: the row that used to be row 3
# becomes row 2, the row that used to be row 4 becomes row 3, and the row
# that used to be row 2 becomes row 4 (A "new week" pushed the oldest of
# the three rolling rows to the bottom of the rotation).
#
# Capture the original values first (columns D and J..Q are the ones that
# actually differ between the three rows; the rest are identical across
# rows 2-4 in this sheet).

$orig = @{}
foreach ($r in 2,3,4) {
    $orig[$r] = @{
        D = $ws.Cells.Item($r, 4).Value2
        J = $ws.Cells.Item($r, 10).Value2
        K = $ws.Cells.Item($r, 11).Value2
        L = $ws.Cells.Item($r, 12).Value2
        M = $ws.Cells.Item($r, 13).Value2
        N = $ws.Cells.Item($r, 14).Value2
        O = $ws.Cells.Item($r, 15).Value2
        P = $ws.Cells.Item($r, 16).Value2
        Q = $ws.Cells.Item($r, 17).Value2
    }
}

# mapping: new row <- old row
$mapping = @{ 2 = 3; 3 = 4; 4 = 2 }

foreach ($newRow in 2,3,4) {
    $src = $orig[$mapping[$newRow]]

    $ws.Cells.Item($newRow, 4).Value2  = $src.D
    $ws.Cells.Item($newRow, 10).Value2 = $src.J
    $ws.Cells.Item($newRow, 11).Value2 = $src.K
    $ws.Cells.Item($newRow, 12).Value2 = $src.L
    $ws.Cells.Item($newRow, 13).Value2 = $src.M
    $ws.Cells.Item($newRow, 14).Value2 = $src.N
    $ws.Cells.Item($newRow, 15).Value2 = $src.O
    $ws.Cells.Item($newRow, 16).Value2 = $src.P
    $ws.Cells.Item($newRow, 17).Value2 = $src.Q
}
